$d = $word.ActiveDocument

# --- 1) Merge the lone leading-space runs into the following text run. ---
# Each of these bodies starts with a standalone " " run followed by a run
# that begins with the given word; "finding & replacing" the space plus the
# first word collapses the two runs into one (same net visible text, just
# re-saved as a single run).
$d.Content.Find.Execute(" Discussed", $true, $false, $false, $false, $false, $true, 1, $false, " Discussed", 2) | Out-Null
$d.Content.Find.Execute(" Changing SKUs", $true, $false, $false, $false, $false, $true, 1, $false, " Changing SKUs", 2) | Out-Null
$d.Content.Find.Execute(" Upon investigation it is determined", $true, $false, $false, $false, $false, $true, 1, $false, " Upon investigation it is determined", 2) | Out-Null
$d.Content.Find.Execute(" Carlo mentioned", $true, $false, $false, $false, $false, $true, 1, $false, " Carlo mentioned", 2) | Out-Null
$d.Content.Find.Execute(" It has been discovered", $true, $false, $false, $false, $false, $true, 1, $false, " It has been discovered", 2) | Out-Null
$d.Content.Find.Execute(" Upon investigation into this process", $true, $false, $false, $false, $false, $true, 1, $false, " Upon investigation into this process", 2) | Out-Null

# --- 2) In the "Return Trade Ins" paragraph, fold the lone-space run between
#        "...stay separate?" and "If we were able..." back into running text. ---
$d.Content.Find.Execute("stay separate? If we were able", $true, $false, $false, $false, $false, $true, 1, $false, "stay separate? If we were able", 2) | Out-Null

# --- 3) "Cashout Revamp" paragraph: append the new trailing sentences. ---
$r = $d.Content
$r.Find.Execute("wouldn’t match.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter(" Do they want to track back to 8/01? If yes")
$r.Collapse(0)
$r.InsertAfter(",")
$r.Collapse(0)
$r.InsertAfter(" then need their printouts for each location back to then.")

# --- 4) "Overall Clean Up and Bug Fix" paragraph: split "more smoothly." so
#        the trailing "ore smoothly." lands in a new run after the _GoBack
#        bookmark, matching the recorded last-edit position. ---
$r2 = $d.Content
$r2.Find.Execute("run more smoothly.", $true, $false, $false, $false, $false, $true, 1, $false, "run m", 2) | Out-Null
$r2.Collapse(0)
$r2.InsertAfter("ore smoothly.")
